# Small preparations for vk parser adding
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill in the VK contact that was previously blank, and normalize
# the "last date" timestamp down to a date-only value.
$ws.Range("E2").Value = "https://vk.com/ff_mgu"
$ws.Range("G2").Value = 46070

# Row 3: normalize the "last date" timestamp down to a date-only value.
$ws.Range("G3").Value = 45910

# New row 4: add the optics/spectroscopy/nanosystems department.
$ws.Range("A4").Value = "Оптики, спектроскопии и физики наносистем"
$ws.Range("B4").Value = "Экспериментальной и теоретической физики"
$ws.Range("C4").Value = "А"
# D4 and F4 stay empty (Сайт / Телеграмм unknown) - touch formatting only
# so the cells exist in the sheet, matching the blank siblings above them.
$ws.Range("D4").Font.Bold = $false
$ws.Range("E4").Value = "http://vk.com/club215377281"
$ws.Range("F4").Font.Bold = $false
$ws.Range("G4").Value = 46056
$ws.Range("G4").NumberFormat = $ws.Range("G3").NumberFormat
